$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$ws.Range("C7").Value = 9
$ws.Range("C11").Value = 10

$ws.Range("C11").Select()
